$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

$b14 = $ws.Cells.Item(14, 2)
$b14.VerticalAlignment = -4108
$b14.Font.Bold = $true
$b14.Font.Size = 12
$b14.Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A14:B17").Select() | Out-Null
